# Purchase Orders sheet rework: new columns (Quantity/Unit Price/Total split
# out of the old "Ordered Quantity" column) and refreshed order rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column widths -------------------------------------------------------
# OOXML <col> width = COM ColumnWidth + 5/6 (Calibri 11 default digit width),
# so subtract 5/6 from every target width below.
$ws.Columns.Item(1).ColumnWidth = 19.166666666666668   # A: 20
$ws.Columns.Item(2).ColumnWidth = 29.166666666666668   # B: 30
$ws.Range("C:E").ColumnWidth     = 14.166666666666666  # C,D,E: 15
$ws.Columns.Item(6).ColumnWidth = 24.166666666666668   # F: 25
$ws.Range("G:I").ColumnWidth     = 14.166666666666666  # G,H,I: 15

# ---- Header row ------------------------------------------------------------
$ws.Range("A1").Value = "Order ID"
$ws.Range("B1").Value = "Product"
$ws.Range("C1").Value = "Quantity"
$ws.Range("D1").Value = "Unit Price"
$ws.Range("E1").Value = "Total"
$ws.Range("F1").Value = "Supplier"
$ws.Range("G1").Value = "Order Date"
$ws.Range("H1").Value = "Estimated Arrival"
$ws.Range("I1").Value = "Payment Status"

# ---- Data rows ---------------------------------------------------------
# Row 2
$ws.Range("A2").Value = "'""67ccdf22f25735bf5559c333"""
$ws.Range("B2").Value = "Industrial Gloves"
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = "'$5.00"
$ws.Range("E2").Value = "'$50.00"
$ws.Range("F2").Value = "Tech Supply Co."
$ws.Range("G2").Value = "'2025-03-08"
$ws.Range("H2").Value = "'2025-03-15"
$ws.Range("I2").Value = "Pending"

# Row 3
$ws.Range("A3").Value = "'""67ccdf22f25735bf5559c333"""
$ws.Range("B3").Value = "Safety Glasses"
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = "'$10.00"
$ws.Range("E3").Value = "'$50.00"
$ws.Range("F3").Value = "Tech Supply Co."
$ws.Range("G3").Value = "'2025-03-08"
$ws.Range("H3").Value = "'2025-03-15"
$ws.Range("I3").Value = "Pending"

# Row 4
$ws.Range("A4").Value = "'""67ccdf22f25735bf5559c336"""
$ws.Range("B4").Value = "Hard Hat"
$ws.Range("C4").Value = 15
$ws.Range("D4").Value = "'$20.00"
$ws.Range("E4").Value = "'$300.00"
$ws.Range("F4").Value = "Manufacture Direct"
$ws.Range("G4").Value = "'2025-03-08"
$ws.Range("H4").Value = "'2025-03-15"
$ws.Range("I4").Value = "Pending"

# Row 5
$ws.Range("A5").Value = "'""67ccdf22f25735bf5559c338"""
$ws.Range("B5").Value = "Safety Glasses"
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = "'$10.00"
$ws.Range("E5").Value = "'$100.00"
$ws.Range("F5").Value = "Tech Supply Co."
$ws.Range("G5").Value = "'2025-03-08"
$ws.Range("H5").Value = "'2025-03-15"
$ws.Range("I5").Value = "Pending"

# The leading apostrophes above force text entry for number/date/currency
# -looking values (matches the source file, where every non-Quantity cell
# is stored as shared-string text). Clear the resulting "quote prefix" cell
# style so the cells end up styleless, same as the rest of the sheet.
$ws.Range("A1:I5").Style = "Normal"
